# edit.ps1 - apply "last commit for a while" changes to Project Update.docx
#
# Strategy: Word's Range.InsertXML (Flat-OPC "pkg:package"-wrapped w:document
# fragment) lets us replace a paragraph's *interior* (the range between its
# start and its paragraph mark) with an exact sequence of <w:r> runs, without
# touching the paragraph's own <w:pPr>. We use this for every paragraph edit
# so the resulting run layout matches the target precisely, and also to
# append the four brand-new list paragraphs at the end of the list (right
# before </w:body> / <w:sectPr>).
#
# NOTE: this runtime's PowerShell seems to mis-bind a parenthesized /
# multi-line expression passed as a second positional argument when the
# first positional argument is a COM object, so every argument is first
# computed into its own plain variable before being passed into a function.

$d = $word.ActiveDocument

$LQ = [char]0x201C   # left curly quote “
$RQ = [char]0x201D   # right curly quote ”

$flatOpcHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$flatOpcFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphRuns($paragraph, [string]$innerRunsXml) {
    # Replace everything inside the paragraph (but not the paragraph mark
    # itself, and not its <w:pPr>) with the given raw <w:r>...</w:r> xml.
    $pStart = $paragraph.Range.Start
    $pEnd = $paragraph.Range.End
    $range = $d.Range($pStart, $pEnd - 1)
    $bodyXml = "<w:p>" + $innerRunsXml + "</w:p>"
    $xml = $flatOpcHeader + $bodyXml + $flatOpcFooter
    $range.InsertXML($xml)
}

function Add-ParagraphsAfter($paragraph, [string]$paragraphsXml) {
    # Insert one or more whole <w:p>...</w:p> paragraphs right after the
    # given paragraph.
    $pEnd = $paragraph.Range.End
    $pos = $d.Range($pEnd, $pEnd)
    $xml = $flatOpcHeader + $paragraphsXml + $flatOpcFooter
    $pos.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Paragraph 1: drop curly quotes around "Covid Data Normalization"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Xml = '<w:r><w:t>Renamed covid_explanatory_analysis file to Covid Data Normalization file</w:t></w:r>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>'
Set-ParagraphRuns $p1 $p1Xml

# ---------------------------------------------------------------------
# Paragraph 2: drop quotes around "Split Tables"; keep curly quotes around
# "Split By Continent" but as their own separate runs.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2Xml = '<w:r><w:t xml:space="preserve">Renamed Split Tables jupyter notebook file to </w:t></w:r>' +
         "<w:r><w:t>$LQ</w:t></w:r>" +
         '<w:r><w:t>Split By Continent</w:t></w:r>' +
         "<w:r><w:t>$RQ</w:t></w:r>"
Set-ParagraphRuns $p2 $p2Xml

# ---------------------------------------------------------------------
# Paragraph 3: merge all runs into one, dropping all curly quotes.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3Xml = '<w:r><w:t>In Covid Data Normalization file split the covid_deaths table into a table that appears once per location, table with columns that are about death (e.g. total_deaths) and a table with the other columns and placed them in dataframes named country, death, other</w:t></w:r>'
Set-ParagraphRuns $p3 $p3Xml

# ---------------------------------------------------------------------
# Paragraph 4: keep two runs, dropping all curly quotes.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4Xml = '<w:r><w:t xml:space="preserve">Placed country data in country.csv, death data in death.csv, other data into </w:t></w:r>' +
         '<w:r><w:t>cases.csv in DB Project folder (should move it out to main file).</w:t></w:r>'
Set-ParagraphRuns $p4 $p4Xml

# ---------------------------------------------------------------------
# Append four brand-new list paragraphs after paragraph 4.
# ---------------------------------------------------------------------
$pPrXml = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>'

$newP5 = '<w:p>' + $pPrXml +
             '<w:r><w:t>Moved DB Project folder (a folder containing the csv files containing all the rows for that database with all tables in 3</w:t></w:r>' +
             '<w:r><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>rd</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> normal form) to the main folder.</w:t></w:r>' +
         '</w:p>'

$newP6 = '<w:p>' + $pPrXml +
             '<w:r><w:t>Made some adjustments to the Covid Data Normalization file.</w:t></w:r>' +
         '</w:p>'

$newP7 = '<w:p>' + $pPrXml +
             '<w:r><w:t>Will only focus on covid death for now.</w:t></w:r>' +
         '</w:p>'

$newP8 = '<w:p>' + $pPrXml +
             '<w:r><w:t>Only performing some SQL queries and adding finishing touches to my tableau</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
             '<w:r><w:t>dashboard and the project will be complete.</w:t></w:r>' +
         '</w:p>'

$newParagraphsXml = $newP5 + $newP6 + $newP7 + $newP8
Add-ParagraphsAfter $p4 $newParagraphsXml
